$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set column D (Price) cells to Text format first so numeric-looking
# strings (e.g. "0.590", "316.71") keep their exact text representation
# instead of being auto-parsed/normalized as numbers by Excel.
$priceRows = @(2,3,5,6,7,9,10,11,12,14,15,16,17,18,19,20,21,22,23,24,25,26,31,32,33,34,36,37,40,41,43,45,46,48,50,51)
foreach ($r in $priceRows) {
    $ws.Cells.Item($r, 4).NumberFormat = "@"
}

# Apply the updated Price (D) and Volume(1h) (E) values row by row.
$ws.Range("D2").Value = "42.844.64"
$ws.Range("E2").Value = "  -5.07%  "
$ws.Range("D3").Value = "2.208.25"
$ws.Range("E3").Value = "  -6.49%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "316.71"
$ws.Range("E5").Value = "  +1.54%  "
$ws.Range("D6").Value = "99.45"
$ws.Range("E6").Value = "  -8.21%  "
$ws.Range("D7").Value = "0.590"
$ws.Range("E7").Value = "  -6.49%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("D9").Value = "0.562"
$ws.Range("E9").Value = "  -7.83%  "
$ws.Range("D10").Value = "36.85"
$ws.Range("E10").Value = "  -9.80%  "
$ws.Range("D11").Value = "53.97"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Value = "0.0829"
$ws.Range("E12").Value = "  -9.53%  "
$ws.Range("E13").Value = "  -8.85%  "
$ws.Range("D14").Value = "0.107"
$ws.Range("E14").Value = "  -2.35%  "
$ws.Range("D15").Value = "0.859"
$ws.Range("D16").Value = "2.547.31"
$ws.Range("E16").Value = "  -6.49%  "
$ws.Range("D17").Value = "14.22"
$ws.Range("E17").Value = "  -6.45%  "
$ws.Range("D18").Value = "2.212.43"
$ws.Range("E18").Value = "  -6.28%  "
$ws.Range("D19").Value = "42.757.52"
$ws.Range("E19").Value = "  -5.20%  "
$ws.Range("D20").Value = "14.62"
$ws.Range("E20").Value = "  +1.29%  "
$ws.Range("D21").Value = "0.0₃0961"
$ws.Range("E21").Value = "  -9.44%  "
$ws.Range("D22").Value = "6.42"
$ws.Range("E22").Value = "  -10.76%  "
$ws.Range("D23").Value = "65.34"
$ws.Range("E23").Value = "  -10.61%  "
$ws.Range("D24").Value = "3.15"
$ws.Range("E24").Value = "  -10.15%  "
$ws.Range("D25").Value = "236.14"
$ws.Range("E25").Value = "  -8.83%  "
$ws.Range("D26").Value = "2.14"
$ws.Range("E26").Value = "  -7.73%  "
$ws.Range("E27").Value = "  +0.13%  "
$ws.Range("E28").Value = "  -9.47%  "
$ws.Range("E29").Value = "  -4.37%  "
$ws.Range("E30").Value = "  -12.32%  "
$ws.Range("D31").Value = "0.0893"
$ws.Range("E31").Value = "  -7.71%  "
$ws.Range("D32").Value = "20.56"
$ws.Range("E32").Value = "  -8.06%  "
$ws.Range("D33").Value = "34.32"
$ws.Range("E33").Value = "  -7.58%  "
$ws.Range("D34").Value = "156.15"
$ws.Range("E34").Value = "  -7.01%  "
$ws.Range("E35").Value = "  -6.77%  "
$ws.Range("D36").Value = "3.15"
$ws.Range("E36").Value = "  +9.75%  "
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  +8.63%  "
$ws.Range("E38").Value = "  -6.55%  "
$ws.Range("E39").Value = "  -8.32%  "
$ws.Range("D40").Value = "4.43"
$ws.Range("E40").Value = "  -5.44%  "
$ws.Range("D41").Value = "3.86"
$ws.Range("E41").Value = "  -2.67%  "
$ws.Range("E42").Value = "  -7.92%  "
$ws.Range("D43").Value = "1.888.90"
$ws.Range("E43").Value = "  +1.39%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("D45").Value = "12.45"
$ws.Range("E45").Value = "  -3.03%  "
$ws.Range("D46").Value = "88.02"
$ws.Range("E46").Value = "  -11.96%  "
$ws.Range("E47").Value = "  -9.05%  "
$ws.Range("D48").Value = "61.23"
$ws.Range("E48").Value = "  -11.83%  "
$ws.Range("E49").Value = "  -5.54%  "
$ws.Range("D50").Value = "76.36"
$ws.Range("E50").Value = "  -8.06%  "
$ws.Range("D51").Value = "102.47"
$ws.Range("E51").Value = "  -7.05%  "
